$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4000
$ws.Range("J74").Value = 5000
$ws.Range("L74").Value = 5000
$ws.Range("N74").Value = -6872
$ws.Range("H77").Value = 4000
$ws.Range("J77").Value = 5000
$ws.Range("L77").Value = 25000
$ws.Range("N77").Value = -34360

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 44999.5
$ws.Range("J24").Value = 44999.5
$ws.Range("L24").Value = 44999.5
$ws.Range("N24").Value = -45747.5
$ws.Range("H45").Value = 3662.875
$ws.Range("I45").Value = 3073.2727
$ws.Range("J45").Value = 4960
$ws.Range("K45").Value = 3073.2727
$ws.Range("L45").Value = 4960
$ws.Range("M45").Value = -2696.2727
$ws.Range("N45").Value = -5714
$ws.Range("H56").Value = 20250
$ws.Range("J56").Value = 20500
$ws.Range("L56").Value = 20500
$ws.Range("N56").Value = -21984
$ws.Range("H61").Value = 5884.65
$ws.Range("I61").Value = 2030.6666
$ws.Range("J61").Value = 11665.625
$ws.Range("K61").Value = 2030.6666
$ws.Range("L61").Value = 11665.625
$ws.Range("M61").Value = -1818.6666
$ws.Range("N61").Value = -12089.625
$ws.Range("H63").Value = 2166.1904
$ws.Range("I63").Value = 2194.0557
$ws.Range("J63").Value = 1999
$ws.Range("K63").Value = 2194.0557
$ws.Range("L63").Value = 1999
$ws.Range("M63").Value = -1508.0557
$ws.Range("N63").Value = -3371
$ws.Range("H66").Value = 2166.1904
$ws.Range("I66").Value = 2194.0557
$ws.Range("J66").Value = 1999
$ws.Range("K66").Value = 10970.2785
$ws.Range("L66").Value = 9995
$ws.Range("M66").Value = -7538.2785
$ws.Range("N66").Value = -16859
$ws.Range("H74").Value = 5749.489
$ws.Range("I74").Value = 985.8570999999999
$ws.Range("J74").Value = 22422.2
$ws.Range("K74").Value = 985.8570999999999
$ws.Range("L74").Value = 22422.2
$ws.Range("M74").Value = -111.8570999999999
$ws.Range("N74").Value = -24170.2
$ws.Range("H77").Value = 5749.489
$ws.Range("I77").Value = 985.8570999999999
$ws.Range("J77").Value = 22422.2
$ws.Range("K77").Value = 4929.2855
$ws.Range("L77").Value = 112111
$ws.Range("M77").Value = -561.2855
$ws.Range("N77").Value = -120847
$ws.Range("H88").Value = 2174.6667
$ws.Range("J88").Value = 1787
$ws.Range("L88").Value = 1787
$ws.Range("N88").Value = -2599
$ws.Range("H91").Value = 2174.6667
$ws.Range("J91").Value = 1787
$ws.Range("L91").Value = 1787
$ws.Range("N91").Value = -4595
$ws.Range("H100").Value = 44999.5
$ws.Range("J100").Value = 44999.5
$ws.Range("L100").Value = 44999.5
$ws.Range("N100").Value = -47163.5
$ws.Range("H132").Value = 2785146.8
$ws.Range("I132").Value = 4171735.2
$ws.Range("K132").Value = 12515205.6
$ws.Range("M132").Value = -12512675.6
$ws.Range("H133").Value = 64475.832
$ws.Range("J133").Value = 64475.832
$ws.Range("L133").Value = 64475.832
$ws.Range("N133").Value = -69535.83199999999
$ws.Range("H136").Value = 5884.65
$ws.Range("I136").Value = 2030.6666
$ws.Range("J136").Value = 11665.625
$ws.Range("K136").Value = 6091.9998
$ws.Range("L136").Value = 34996.875
$ws.Range("M136").Value = -3541.9998
$ws.Range("N136").Value = -40096.875
$ws.Range("H139").Value = 99999
$ws.Range("J139").Value = 99999
$ws.Range("L139").Value = 99999
$ws.Range("N139").Value = -110279

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2777.7368
$ws.Range("I86").Value = 2675.2307
$ws.Range("J86").Value = 2999.8333
$ws.Range("K86").Value = 2675.2307
$ws.Range("L86").Value = 2999.8333
$ws.Range("M86").Value = -1552.2307
$ws.Range("N86").Value = -5245.8333
$ws.Range("H89").Value = 2777.7368
$ws.Range("I89").Value = 2675.2307
$ws.Range("J89").Value = 2999.8333
$ws.Range("K89").Value = 13376.1535
$ws.Range("L89").Value = 14999.1665
$ws.Range("M89").Value = -7760.1535
$ws.Range("N89").Value = -26231.1665

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2144.7454
$ws.Range("I31").Value = 2066.625
$ws.Range("J31").Value = 2253.4348
$ws.Range("K31").Value = 2066.625
$ws.Range("L31").Value = 2253.4348
$ws.Range("M31").Value = -1771.625
$ws.Range("N31").Value = -2843.4348
$ws.Range("H34").Value = 2144.7454
$ws.Range("I34").Value = 2066.625
$ws.Range("J34").Value = 2253.4348
$ws.Range("K34").Value = 2066.625
$ws.Range("L34").Value = 2253.4348
$ws.Range("M34").Value = -1864.625
$ws.Range("N34").Value = -2657.4348
$ws.Range("H62").Value = 3468.3
$ws.Range("I62").Value = 3740.2856
$ws.Range("K62").Value = 3740.2856
$ws.Range("M62").Value = -3116.2856
$ws.Range("H65").Value = 3468.3
$ws.Range("I65").Value = 3740.2856
$ws.Range("K65").Value = 18701.428
$ws.Range("M65").Value = -15581.428

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 2841.3333
$ws.Range("I99").Value = 1409.6
$ws.Range("J99").Value = 10000
$ws.Range("K99").Value = 4228.799999999999
$ws.Range("L99").Value = 30000
$ws.Range("M99").Value = -1982.799999999999
$ws.Range("N99").Value = -34492
$ws.Range("H131").Value = 2843.7354
$ws.Range("J131").Value = 2900.2122
$ws.Range("L131").Value = 8700.6366
$ws.Range("N131").Value = -18780.6366
$ws.Range("H132").Value = 740.2727
$ws.Range("I132").Value = 599
$ws.Range("J132").Value = 987.5
$ws.Range("K132").Value = 5391
$ws.Range("L132").Value = 8887.5
$ws.Range("M132").Value = -2861
$ws.Range("N132").Value = -13947.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 884.6486
$ws.Range("I97").Value = 853.6177
$ws.Range("K97").Value = 853.6177
$ws.Range("M97").Value = -357.6177
$ws.Range("H113").Value = 3120.2942
$ws.Range("I113").Value = 2251.7273
$ws.Range("K113").Value = 2251.7273
$ws.Range("M113").Value = -81.72730000000001
$ws.Range("H132").Value = 14020.526
$ws.Range("I132").Value = 18855.72
$ws.Range("J132").Value = 4722.077
$ws.Range("K132").Value = 56567.16
$ws.Range("L132").Value = 14166.231
$ws.Range("M132").Value = -54037.16
$ws.Range("N132").Value = -19226.231

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").Value = $null
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").Value = $null
$ws.Range("H82").Value = 1344.5333
$ws.Range("I82").Value = 1760.1428
$ws.Range("K82").Value = 1760.1428
$ws.Range("M82").Value = -1399.1428
$ws.Range("H85").Value = 1344.5333
$ws.Range("I85").Value = 1760.1428
$ws.Range("K85").Value = 1760.1428
$ws.Range("M85").Value = -512.1428000000001
$ws.Range("H93").Value = 4312.3335
$ws.Range("I93").Value = 1484.25
$ws.Range("J93").Value = 9968.5
$ws.Range("K93").Value = 1484.25
$ws.Range("L93").Value = 9968.5
$ws.Range("M93").Value = -236.25
$ws.Range("N93").Value = -12464.5
$ws.Range("H110").Value = 32500
$ws.Range("J110").Value = 32500
$ws.Range("L110").Value = 32500
$ws.Range("N110").Value = -40680
$ws.Range("H132").Value = 3006.6875
$ws.Range("I132").Value = 2925.25
$ws.Range("J132").Value = 3088.125
$ws.Range("K132").Value = 8775.75
$ws.Range("L132").Value = 9264.375
$ws.Range("M132").Value = -6245.75
$ws.Range("N132").Value = -14324.375
$ws.Range("H133").Value = 67500
$ws.Range("J133").Value = 67500
$ws.Range("L133").Value = 67500
$ws.Range("N133").Value = -72560
$ws.Range("H136").Value = 8555.954
$ws.Range("I136").Value = 4378.273
$ws.Range("J136").Value = 12733.637
$ws.Range("K136").Value = 13134.819
$ws.Range("L136").Value = 38200.911
$ws.Range("M136").Value = -10584.819
$ws.Range("N136").Value = -43300.911

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 140000
$ws.Range("J57").Value = 140000
$ws.Range("L57").Value = 140000
$ws.Range("N57").Value = -141508
$ws.Range("H96").Value = 4239.9
$ws.Range("I96").Value = 4133.3335
$ws.Range("K96").Value = 4133.3335
$ws.Range("M96").Value = -2760.3335
$ws.Range("H100").Value = 1433.4546
$ws.Range("I100").Value = 1596.25
$ws.Range("K100").Value = 3192.5
$ws.Range("M100").Value = -2651.5
